# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# on several rows across multiple sheets, per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 491.18182
$ws.Range("I33").Value = 240.3
$ws.Range("K33").Value = 240.3
$ws.Range("M33").Value = -11.30000000000001

$ws.Range("H40").Value = 4056.6667
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4056.6667
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4056.6667
$ws.Range("N40").Value = -4406.6667
$ws.Range("M40").ClearContents()

$ws.Range("H100").Value = 2959.6
$ws.Range("I100").Value = 2959.6
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2959.6
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2418.6
$ws.Range("N100").ClearContents()

$ws.Range("H125").Value = 2986.2856
$ws.Range("I125").Value = 1484
$ws.Range("K125").Value = 13356
$ws.Range("M125").Value = -10896

$ws.Range("H135").Value = 1322.7778
$ws.Range("I135").Value = 1361
$ws.Range("J135").Value = 1292.2
$ws.Range("K135").Value = 12249
$ws.Range("L135").Value = 11629.8
$ws.Range("M135").Value = -9714
$ws.Range("N135").Value = -16699.8

$ws.Range("H138").Value = 3455.111
$ws.Range("I138").Value = 2732
$ws.Range("J138").Value = 3816.6667
$ws.Range("K138").Value = 8196
$ws.Range("L138").Value = 11450.0001
$ws.Range("M138").Value = -3056
$ws.Range("N138").Value = -21730.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 13171.75
$ws.Range("J96").Value = 13171.75
$ws.Range("L96").Value = 13171.75
$ws.Range("N96").Value = -18663.75

$ws.Range("H131").Value = 71250
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 71250
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 71250
$ws.Range("N131").Value = -81330
$ws.Range("M131").ClearContents()

$ws.Range("H132").Value = 2840.6428
$ws.Range("I132").Value = 2905.75
$ws.Range("K132").Value = 8717.25
$ws.Range("M132").Value = -6187.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 389
$ws.Range("J80").Value = 533.5
$ws.Range("L80").Value = 533.5
$ws.Range("N80").Value = -2529.5

$ws.Range("H83").Value = 389
$ws.Range("J83").Value = 533.5
$ws.Range("L83").Value = 2667.5
$ws.Range("N83").Value = -12651.5

$ws.Range("H86").Value = 6000
$ws.Range("J86").Value = 6750
$ws.Range("L86").Value = 6750
$ws.Range("N86").Value = -8996

$ws.Range("H89").Value = 6000
$ws.Range("J89").Value = 6750
$ws.Range("L89").Value = 33750
$ws.Range("N89").Value = -44982

$ws.Range("H107").Value = 1136.8
$ws.Range("I107").Value = 1172.25
$ws.Range("J107").Value = 995
$ws.Range("K107").Value = 1172.25
$ws.Range("L107").Value = 995
$ws.Range("M107").Value = 747.75
$ws.Range("N107").Value = -4835

$ws.Range("H127").Value = 80000
$ws.Range("J127").Value = 80000
$ws.Range("L127").Value = 80000
$ws.Range("N127").Value = -89920

$ws.Range("H134").Value = 6149.2583
$ws.Range("I134").Value = 6556.6665
$ws.Range("J134").Value = 3399.25
$ws.Range("K134").Value = 19669.9995
$ws.Range("L134").Value = 10197.75
$ws.Range("M134").Value = -17134.9995
$ws.Range("N134").Value = -15267.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 27696.182
$ws.Range("I59").Value = 23522.572
$ws.Range("K59").Value = 23522.572
$ws.Range("M59").Value = -22377.572

$ws.Range("H115").Value = 40290
$ws.Range("J115").Value = 40290
$ws.Range("L115").Value = 40290
$ws.Range("N115").Value = -42640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 264.66666
$ws.Range("I12").Value = 76.333336
$ws.Range("J12").Value = 358.83334
$ws.Range("K12").Value = 229.000008
$ws.Range("L12").Value = 1076.50002
$ws.Range("M12").Value = -56.00000800000001
$ws.Range("N12").Value = -1422.50002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6099.8335
$ws.Range("I70").Value = 5524.75
$ws.Range("K70").Value = 5524.75
$ws.Range("M70").Value = -5254.75

$ws.Range("H73").Value = 6099.8335
$ws.Range("I73").Value = 5524.75
$ws.Range("K73").Value = 5524.75
$ws.Range("M73").Value = -4588.75

$ws.Range("H80").Value = 9018.166999999999
$ws.Range("J80").Value = 9578.333000000001
$ws.Range("L80").Value = 9578.333000000001
$ws.Range("N80").Value = -11574.333

$ws.Range("H83").Value = 9018.166999999999
$ws.Range("J83").Value = 9578.333000000001
$ws.Range("L83").Value = 47891.665
$ws.Range("N83").Value = -57875.665

$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988

$ws.Range("H122").Value = 56744
$ws.Range("I122").Value = 59766.4
$ws.Range("J122").Value = 46669.332
$ws.Range("K122").Value = 179299.2
$ws.Range("L122").Value = 140007.996
$ws.Range("M122").Value = -176849.2
$ws.Range("N122").Value = -144907.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 600
$ws.Range("I7").Value = 200
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 200
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -88
$ws.Range("N7").Value = -1224

$ws.Range("H40").Value = 5174
$ws.Range("I40").Value = 5213.2856
$ws.Range("K40").Value = 5213.2856
$ws.Range("M40").Value = -5077.2856

$ws.Range("H68").Value = 733.3333
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 1000
$ws.Range("N68").Value = -2498

$ws.Range("H71").Value = 733.3333
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 5000
$ws.Range("N71").Value = -12488

$ws.Range("H126").Value = 600
$ws.Range("I126").Value = 200
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 600
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 1870
$ws.Range("N126").Value = -7940

$ws.Range("H132").Value = 21807.111
$ws.Range("I132").Value = 21784
$ws.Range("J132").Value = 21888
$ws.Range("K132").Value = 65352
$ws.Range("L132").Value = 65664
$ws.Range("M132").Value = -62822
$ws.Range("N132").Value = -70724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4000
$ws.Range("J81").Value = 4000
$ws.Range("L81").Value = 8000
$ws.Range("N81").Value = -10122

$ws.Range("H84").Value = 4000
$ws.Range("J84").Value = 4000
$ws.Range("L84").Value = 40000
$ws.Range("N84").Value = -50608

$ws.Range("H135").Value = 49499.875
$ws.Range("J135").Value = 53714.145
$ws.Range("L135").Value = 53714.145
$ws.Range("N135").Value = -63854.145
